$d = $word.ActiveDocument

# --- Fix paragraph 1: "Devair Restani" -> "X"; rejoin split word "pa|ssar"; "rio" -> "assude" ---
$d.Content.Find.Execute("Devair Restani", $true, $false, $false, $false, $false, $true, 1, $false, "X", 2)
$d.Content.Find.Execute("pescar e pa" + "ssar o dia", $true, $false, $false, $false, $false, $true, 1, $false, "pescar e passar o dia", 2)
$d.Content.Find.Execute("aproveitar o dia beira rio.", $true, $false, $false, $false, $false, $true, 1, $false, "aproveitar o dia beira assude.", 2)

# --- Fix paragraph 2: "errada" -> "equivocada" ---
$d.Content.Find.Execute("inserção errada de produtos", $true, $false, $false, $false, $false, $true, 1, $false, "inserção equivocada de produtos", 2)

# --- Fix paragraph 3: "frequencia" -> "frequência" ---
$d.Content.Find.Execute("a frequencia de erros", $true, $false, $false, $false, $false, $true, 1, $false, "a frequência de erros", 2)

# --- Insert three new content paragraphs after paragraph 4 ---
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertAfter("Para a gerência da comanda dos clientes o sistema deverá fornecer opções de inserir produtos nela, privando-os do estoque, para a mesma. Como cada produto possui um preço por unidade, a comanda efetuará a soma de todos os inseridos para a venda destes.")

$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertAfter("O fechamento da comanda será efetuado após o cliente pagar o valor equivalente ao valor total da soma dos produtos que foram inseridos na comanda. ")

$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertAfter("Contudo o sistema deve suprir as necessidades básicas de gerenciamento do pesqueiro X, evitando erros humanos, assim gerando uma maior segurança para com o lucro de seu negócio.")

# --- Move the _GoBack bookmark to sit between "negócio" and the final "." ---
$findRange = $d.Paragraphs.Item(7).Range.Duplicate
$findRange.Find.Execute("negócio.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dotPos = $findRange.End - 1
$bmRange = $d.Range($dotPos, $dotPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Append two empty paragraphs (no run at all) and one paragraph containing a single space ---
$blankXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="420" w:firstLineChars="0"/><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p>'

$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertXML($blankXml)

$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertXML($blankXml)

$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertAfter(" ")
